# Updates for version 3 compatibility
# - Add two new "REVIEW" type annotation rows to the Annotations sheet
#   (migrating the old Reviewers-sheet entries into the Annotations model)
# - Make the Annotations sheet the active/selected sheet, with A16 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Annotations")

# Make this sheet the active one (also flips tabSelected off the
# previously-active sheet and onto this one, and updates activeTab).
$ws.Activate() | Out-Null

# Insert two new blank rows before the existing row 3 (pushing the old
# rows 3 and 4 down to rows 5 and 6).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# New row 3: reviewer "Joe" migrated from the Reviewers sheet.
$ws.Cells.Item(3, 1).Value = "SPDXRef-DOCUMENT"
$ws.Cells.Item(3, 2).Value = "This is just an example.  Some of the non-standard licenses look like they are actually BSD 3 clause licenses"
$ws.Cells.Item(3, 3).Value = "2010-02-10T00:00:00Z"
$ws.Cells.Item(3, 4).Value = "Person: Joe Reviewer"
$ws.Cells.Item(3, 5).Value = "REVIEW"

# New row 4: reviewer "Suzanne" migrated from the Reviewers sheet.
$ws.Cells.Item(4, 1).Value = "SPDXRef-DOCUMENT"
$ws.Cells.Item(4, 2).Value = "Another example reviewer."
$ws.Cells.Item(4, 3).Value = "2011-03-13T00:00:00Z"
$ws.Cells.Item(4, 4).Value = "Person: Suzanne Reviewer"
$ws.Cells.Item(4, 5).Value = "REVIEW"

# Leave the selection on A16, matching the saved view state.
$ws.Range("A16").Select() | Out-Null
